$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Clear the text from A1 (header becomes blank, matching diff)
$ws.Range("A1").ClearContents()

# 2. Remove header formatting (bold font, thin border, center/top alignment)
#    so header cells revert to the default style (style 0).
$ws.Range("A1:AP1").ClearFormats()

# 3. Update numeric cell values for rows 3-7 per corrected data cleaning
# Row 3
$ws.Range("B3").Value = 3
$ws.Range("F3").Value = 3
$ws.Range("H3").Value = 29
$ws.Range("I3").Value = 10
$ws.Range("K3").Value = 30
$ws.Range("L3").Value = 15
$ws.Range("N3").Value = 10
$ws.Range("R3").Value = 3
$ws.Range("S3").Value = 27
$ws.Range("T3").Value = 14
$ws.Range("U3").Value = 6
$ws.Range("W3").Value = 12
$ws.Range("X3").Value = 4
$ws.Range("Z3").Value = 13
$ws.Range("AB3").Value = 23
$ws.Range("AM3").Value = 0
$ws.Range("AP3").Value = 8

# Row 4
$ws.Range("B4").Value = 5
$ws.Range("F4").Value = 4
$ws.Range("H4").Value = 112
$ws.Range("I4").Value = 22
$ws.Range("K4").Value = 84
$ws.Range("L4").Value = 24
$ws.Range("N4").Value = 11
$ws.Range("R4").Value = 4
$ws.Range("S4").Value = 64
$ws.Range("T4").Value = 24
$ws.Range("U4").Value = 10
$ws.Range("W4").Value = 18
$ws.Range("X4").Value = 6
$ws.Range("Z4").Value = 21
$ws.Range("AB4").Value = 61
$ws.Range("AM4").Value = 1
$ws.Range("AP4").Value = 9

# Row 5
$ws.Range("B5").Value = 3303.59
$ws.Range("F5").Value = 2944.79
$ws.Range("H5").Value = 46586
$ws.Range("I5").Value = 11488.1
$ws.Range("K5").Value = 35412.62
$ws.Range("L5").Value = 10778.18
$ws.Range("N5").Value = 5939.57
$ws.Range("R5").Value = 2944.79
$ws.Range("S5").Value = 26936.95
$ws.Range("T5").Value = 12546.44
$ws.Range("U5").Value = 5080.62
$ws.Range("W5").Value = 8258.83
$ws.Range("X5").Value = 4054.68
$ws.Range("Z5").Value = 10455.85
$ws.Range("AB5").Value = 24870.33
$ws.Range("AM5").Value = 316.76
$ws.Range("AP5").Value = 5489.07

# Row 6
$ws.Range("B6").Value = 2.74
$ws.Range("D6").Value = 0.17
$ws.Range("F6").Value = 2.45
$ws.Range("G6").Value = 0.97
$ws.Range("H6").Value = 38.7
$ws.Range("I6").Value = 9.539999999999999
$ws.Range("J6").Value = 0.53
$ws.Range("K6").Value = 29.42
$ws.Range("L6").Value = 8.949999999999999
$ws.Range("M6").Value = 4.93
$ws.Range("N6").Value = 4.93
$ws.Range("P6").Value = 0.97
$ws.Range("R6").Value = 2.45
$ws.Range("S6").Value = 22.38
$ws.Range("T6").Value = 10.42
$ws.Range("U6").Value = 4.22
$ws.Range("V6").Value = 3.4
$ws.Range("W6").Value = 6.86
$ws.Range("X6").Value = 3.37
$ws.Range("Z6").Value = 8.69
$ws.Range("AA6").Value = 1.44
$ws.Range("AB6").Value = 20.66
$ws.Range("AE6").Value = 2.45
$ws.Range("AF6").Value = 1.25
$ws.Range("AH6").Value = 0.39
$ws.Range("AI6").Value = 2.23
$ws.Range("AK6").Value = 0.97
$ws.Range("AM6").Value = 0.26
$ws.Range("AP6").Value = 4.56

# Row 7
$ws.Range("B7").Value = 660.72
$ws.Range("F7").Value = 736.2
$ws.Range("H7").Value = 415.95
$ws.Range("I7").Value = 522.1900000000001
$ws.Range("K7").Value = 421.58
$ws.Range("L7").Value = 449.09
$ws.Range("N7").Value = 539.96
$ws.Range("R7").Value = 736.2
$ws.Range("S7").Value = 420.89
$ws.Range("T7").Value = 522.77
$ws.Range("U7").Value = 508.06
$ws.Range("W7").Value = 458.82
$ws.Range("X7").Value = 675.78
$ws.Range("Z7").Value = 497.9
$ws.Range("AB7").Value = 407.71
$ws.Range("AM7").Value = 316.76
$ws.Range("AP7").Value = 609.9

# 4. Delete now-unused trailing rows 9-12 (data only went through row 8)
$ws.Range("A9:AP12").Delete()
